$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Dkk1/Kremen2 signalling recomputed with updated TPM values; ---
# --- the target cluster for this row also changed from FAPs to ECs.    ---
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Dkk1"
$ws.Range("C2").Value = "Kremen2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03046
$ws.Range("H2").Value = 0.09138
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.52656
$ws.Range("N2").Value = 1.57968
$ws.Range("O2").Value = 0.9686255056421601
$ws.Range("P2").Value = 0.9686255056421602
$ws.Range("Q2").Value = 0.0160390176
$ws.Range("R2").Value = 0.1443511584
$ws.Range("S2").Value = 0.9686255056421601
$ws.Range("T2").Value = 0.9686255056421602

# --- Row 3 (new): same Dkk1/Kremen2 pair, target cluster "FAPs" (this ---
# --- is the original row 2's target/values, now on its own row).      ---
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Dkk1"
$ws.Range("C3").Value = "Kremen2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03046
$ws.Range("H3").Value = 0.09138
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01705566666666667
$ws.Range("N3").Value = 0.051167
$ws.Range("O3").Value = 0.03137449435783982
$ws.Range("P3").Value = 0.03137449435783982
$ws.Range("Q3").Value = 0.0005195156066666667
$ws.Range("R3").Value = 0.00467564046
$ws.Range("S3").Value = 0.03137449435783982
$ws.Range("T3").Value = 0.03137449435783982
